$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Bolton & Kacperczyk citation year 2022 -> 2023
$ws.Range("B3").Value = "Patrick Bolton, `nMarcin Kacperczyk;`n2023;`nGlobal Pricing of Carbon-Transition Risk"

# Fill in the RQs/Objectives/Aim for the portfolio paper (row 4), replacing the "To fill" placeholder
$ws.Range("I4").Value = "Compute the BMG risk factor on `nan investment universe of <2000 stocks through the traditional Fama-French approach. Analyse this factor across different contexts, including temporally.`nThen, based on this analytical framework, produce optimisations to portfolio management. Including with regard to minimum variance portfolios and index reweighting. This will further the theme of not bearing unrewarded risk by integrating the further dimension of carbon risk."

# Add the "Applicable to my content" notes for the portfolio paper (row 4)
$ws.Range("Q4").Value = "Portfolio techniques portion is not `nrelevant to my work.`nThe considerations from the evaluation of the risk-factor approach will be broadly applicable to my work."

# These two cells should wrap their text like the rest of the row
$ws.Range("I4").WrapText = $true
$ws.Range("Q4").WrapText = $true

# Move the active selection/view to reflect where the author ended up editing
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Application.ActiveWindow.ScrollColumn = 12
$ws.Range("Q4").Select() | Out-Null
